# Commit: use province as custom individual fixed effect on CHARLS dataset;
# adds a new quadratic income term (AVGINDIINCOME_EARN2 -> RC5/RC6/RC7 PCA
# components) to the CHARLS PCA-loading tables on both worksheets, and drops
# AVGBMI (sheet "OUTP1M_RATIO") / SMOKEEVER_RATIO (sheet "CHRONIC_RATIO").
# Rewritten below as a full refresh of each sheet's used range so every
# row label / numeric loading lines up with the regenerated PCA output.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OUTP1M_RATIO")
$ws2 = $wb.Worksheets.Item("CHRONIC_RATIO")

# --- Sheet "OUTP1M_RATIO": header row + 16 variable rows, columns A:G ---
$ws1.Cells.Item(1,2).Value = 'RC2'
$ws1.Cells.Item(1,3).Value = 'RC1'
$ws1.Cells.Item(1,4).Value = 'RC3'
$ws1.Cells.Item(1,5).Value = 'RC6'
$ws1.Cells.Item(1,6).Value = 'RC5'
$ws1.Cells.Item(1,7).Value = 'RC4'
$ws1.Cells.Item(2,1).Value = 'AVGAGE'
$ws1.Cells.Item(2,2).Value = 0.88718668615851
$ws1.Cells.Item(2,3).Value = -0.144144988547339
$ws1.Cells.Item(2,4).Value = -0.0571651818782883
$ws1.Cells.Item(2,5).Value = -0.148681919740619
$ws1.Cells.Item(2,6).Value = 0.0460349452946635
$ws1.Cells.Item(2,7).Value = -0.0286893866487868
$ws1.Cells.Item(3,1).Value = 'MALE_RATIO'
$ws1.Cells.Item(3,2).Value = -0.163093369845616
$ws1.Cells.Item(3,3).Value = 0.164768396646038
$ws1.Cells.Item(3,4).Value = 0.0253606416552824
$ws1.Cells.Item(3,5).Value = 0.245646139327318
$ws1.Cells.Item(3,6).Value = -0.102046735919805
$ws1.Cells.Item(3,7).Value = 0.612912368676942
$ws1.Cells.Item(4,1).Value = 'MARITAL_AVELEN'
$ws1.Cells.Item(4,2).Value = 0.930741143462427
$ws1.Cells.Item(4,3).Value = 0.120078329350535
$ws1.Cells.Item(4,4).Value = -0.0526608527245032
$ws1.Cells.Item(4,5).Value = 0.113640708187669
$ws1.Cells.Item(4,6).Value = 0.0609799315767796
$ws1.Cells.Item(4,7).Value = -0.0335123295361416
$ws1.Cells.Item(5,1).Value = 'DRINK1Y_RATIO'
$ws1.Cells.Item(5,2).Value = 0.0794150467068352
$ws1.Cells.Item(5,3).Value = -0.0709735050028291
$ws1.Cells.Item(5,4).Value = 0.0604687726386141
$ws1.Cells.Item(5,5).Value = -0.0575069154848589
$ws1.Cells.Item(5,6).Value = 0.00674016334467371
$ws1.Cells.Item(5,7).Value = 0.763223249645825
$ws1.Cells.Item(6,1).Value = 'SMOKENOW_RATIO'
$ws1.Cells.Item(6,2).Value = -0.0355656539824122
$ws1.Cells.Item(6,3).Value = 0.0815696890053038
$ws1.Cells.Item(6,4).Value = -0.0906592121399536
$ws1.Cells.Item(6,5).Value = 0.928449895532468
$ws1.Cells.Item(6,6).Value = -0.144927736996376
$ws1.Cells.Item(6,7).Value = 0.0563140753271835
$ws1.Cells.Item(7,1).Value = 'AVGSMOKENUM'
$ws1.Cells.Item(7,2).Value = 0.0479137590231844
$ws1.Cells.Item(7,3).Value = -0.00734072568255442
$ws1.Cells.Item(7,4).Value = -0.0306017081763378
$ws1.Cells.Item(7,5).Value = 0.951871197755935
$ws1.Cells.Item(7,6).Value = -0.0354345643975689
$ws1.Cells.Item(7,7).Value = 0.0802341687121033
$ws1.Cells.Item(8,1).Value = 'AVGHOSP1Y_REALEXP'
$ws1.Cells.Item(8,2).Value = 0.180016096356673
$ws1.Cells.Item(8,3).Value = -0.206459834489431
$ws1.Cells.Item(8,4).Value = -0.0117083748859195
$ws1.Cells.Item(8,5).Value = -0.0648150835681795
$ws1.Cells.Item(8,6).Value = 0.679902686199386
$ws1.Cells.Item(8,7).Value = 0.0323152739056778
$ws1.Cells.Item(9,1).Value = 'AVGOUTP1M_REALEXP'
$ws1.Cells.Item(9,2).Value = 0.0318785630994733
$ws1.Cells.Item(9,3).Value = 0.0936566498714934
$ws1.Cells.Item(9,4).Value = -0.0115961193886837
$ws1.Cells.Item(9,5).Value = 0.0301634741773827
$ws1.Cells.Item(9,6).Value = 0.728047200041228
$ws1.Cells.Item(9,7).Value = -0.0997849498157468
$ws1.Cells.Item(10,1).Value = 'INSURANCE_RATIO'
$ws1.Cells.Item(10,2).Value = 0.00897058717442542
$ws1.Cells.Item(10,3).Value = 0.050609554446579
$ws1.Cells.Item(10,4).Value = 0.957087527339415
$ws1.Cells.Item(10,5).Value = -0.0638350540379625
$ws1.Cells.Item(10,6).Value = -0.0432722050686113
$ws1.Cells.Item(10,7).Value = 0.0945630241064365
$ws1.Cells.Item(11,1).Value = 'INSGOV_RATIO'
$ws1.Cells.Item(11,2).Value = 0.00958883801332884
$ws1.Cells.Item(11,3).Value = 0.143174876751917
$ws1.Cells.Item(11,4).Value = 0.947272538245271
$ws1.Cells.Item(11,5).Value = -0.0555402268699357
$ws1.Cells.Item(11,6).Value = -0.043161672239431
$ws1.Cells.Item(11,7).Value = 0.0874040327412882
$ws1.Cells.Item(12,1).Value = 'AVGEXP1Y_TOTAL'
$ws1.Cells.Item(12,2).Value = 0.0435693952903351
$ws1.Cells.Item(12,3).Value = -0.241081955484789
$ws1.Cells.Item(12,4).Value = -0.071280050443506
$ws1.Cells.Item(12,5).Value = -0.200369904330724
$ws1.Cells.Item(12,6).Value = 0.648817264968937
$ws1.Cells.Item(12,7).Value = 0.103016683051024
$ws1.Cells.Item(13,1).Value = 'TRANSCHILD_RATIO'
$ws1.Cells.Item(13,2).Value = 0.758819220050551
$ws1.Cells.Item(13,3).Value = 0.270101419135645
$ws1.Cells.Item(13,4).Value = 0.124923517222997
$ws1.Cells.Item(13,5).Value = 0.0322916752291758
$ws1.Cells.Item(13,6).Value = 0.167544734891192
$ws1.Cells.Item(13,7).Value = 0.130969425481883
$ws1.Cells.Item(14,1).Value = 'WORK_RATIO'
$ws1.Cells.Item(14,2).Value = -0.069856165992134
$ws1.Cells.Item(14,3).Value = 0.779456164044369
$ws1.Cells.Item(14,4).Value = 0.213636888848949
$ws1.Cells.Item(14,5).Value = 0.109250322441981
$ws1.Cells.Item(14,6).Value = -0.212614379271422
$ws1.Cells.Item(14,7).Value = 0.318114888851527
$ws1.Cells.Item(15,1).Value = 'JOBSTATUS_AGRI_RATIO'
$ws1.Cells.Item(15,2).Value = 0.0672251437783013
$ws1.Cells.Item(15,3).Value = 0.947225263339962
$ws1.Cells.Item(15,4).Value = 0.152475040575057
$ws1.Cells.Item(15,5).Value = 0.064262085817992
$ws1.Cells.Item(15,6).Value = -0.147910484001572
$ws1.Cells.Item(15,7).Value = 0.114396252801658
$ws1.Cells.Item(16,1).Value = 'JOBSTATUS_NAGE_RATIO'
$ws1.Cells.Item(16,2).Value = -0.210027159120502
$ws1.Cells.Item(16,3).Value = -0.779392327791553
$ws1.Cells.Item(16,4).Value = 0.0726500630155932
$ws1.Cells.Item(16,5).Value = 0.0465349959475207
$ws1.Cells.Item(16,6).Value = 0.0146864981041822
$ws1.Cells.Item(16,7).Value = 0.244495733004015
$ws1.Cells.Item(17,1).Value = 'JOBSTATUS_NEWK_RATIO'
$ws1.Cells.Item(17,2).Value = -0.386805586205041
$ws1.Cells.Item(17,3).Value = -0.0524865907743857
$ws1.Cells.Item(17,4).Value = -0.208758746015093
$ws1.Cells.Item(17,5).Value = -0.0356960506666224
$ws1.Cells.Item(17,6).Value = -0.220967106897498
$ws1.Cells.Item(17,7).Value = -0.507550045905446

# --- Sheet "CHRONIC_RATIO": header row + 17 variable rows, columns A:H ---
$ws2.Cells.Item(1,2).Value = 'RC5'
$ws2.Cells.Item(1,3).Value = 'RC1'
$ws2.Cells.Item(1,4).Value = 'RC2'
$ws2.Cells.Item(1,5).Value = 'RC3'
$ws2.Cells.Item(1,6).Value = 'RC4'
$ws2.Cells.Item(1,7).Value = 'RC6'
$ws2.Cells.Item(1,8).Value = 'RC7'
$ws2.Cells.Item(2,1).Value = 'AVGAGE'
$ws2.Cells.Item(2,2).Value = -0.214572231784796
$ws2.Cells.Item(2,3).Value = 0.881108400901217
$ws2.Cells.Item(2,4).Value = 0.067128881153954
$ws2.Cells.Item(2,5).Value = 0.0786309094527856
$ws2.Cells.Item(2,6).Value = -0.0063294957543317
$ws2.Cells.Item(2,7).Value = 0.102377302889804
$ws2.Cells.Item(2,8).Value = -0.0350045235121683
$ws2.Cells.Item(3,1).Value = 'MARITAL_RATIO'
$ws2.Cells.Item(3,2).Value = -0.152296773250094
$ws2.Cells.Item(3,3).Value = -0.43800797031574
$ws2.Cells.Item(3,4).Value = -0.328374822283223
$ws2.Cells.Item(3,5).Value = 0.19569244420252
$ws2.Cells.Item(3,6).Value = 0.394604255628119
$ws2.Cells.Item(3,7).Value = 0.186753371792834
$ws2.Cells.Item(3,8).Value = 0.215523809034415
$ws2.Cells.Item(4,1).Value = 'MARITAL_AVELEN'
$ws2.Cells.Item(4,2).Value = -0.156015839138029
$ws2.Cells.Item(4,3).Value = 0.870614015687656
$ws2.Cells.Item(4,4).Value = -0.0072944306420246
$ws2.Cells.Item(4,5).Value = -0.144258071326179
$ws2.Cells.Item(4,6).Value = 0.0626276126006858
$ws2.Cells.Item(4,7).Value = 0.0762957496917967
$ws2.Cells.Item(4,8).Value = 0.151527414268612
$ws2.Cells.Item(5,1).Value = 'DRINK1Y_RATIO'
$ws2.Cells.Item(5,2).Value = -0.283644410122958
$ws2.Cells.Item(5,3).Value = -0.0657558626303956
$ws2.Cells.Item(5,4).Value = 0.00668946257008363
$ws2.Cells.Item(5,5).Value = 0.0926256933469915
$ws2.Cells.Item(5,6).Value = 0.723740459576101
$ws2.Cells.Item(5,7).Value = -0.0122085196940393
$ws2.Cells.Item(5,8).Value = -0.1154621832565
$ws2.Cells.Item(6,1).Value = 'AVGHOSP1Y_REALEXP'
$ws2.Cells.Item(6,2).Value = 0.104170197571148
$ws2.Cells.Item(6,3).Value = 0.253436966881592
$ws2.Cells.Item(6,4).Value = 0.316469527961999
$ws2.Cells.Item(6,5).Value = 0.166254590818565
$ws2.Cells.Item(6,6).Value = 0.142943174497681
$ws2.Cells.Item(6,7).Value = 0.561661752196503
$ws2.Cells.Item(6,8).Value = -0.0350414706084805
$ws2.Cells.Item(7,1).Value = 'AVGOUTP1M_REALEXP'
$ws2.Cells.Item(7,2).Value = -0.108941152697723
$ws2.Cells.Item(7,3).Value = -0.0277582912661303
$ws2.Cells.Item(7,4).Value = 0.102788158984461
$ws2.Cells.Item(7,5).Value = -0.0953617522070969
$ws2.Cells.Item(7,6).Value = -0.0570145247764505
$ws2.Cells.Item(7,7).Value = 0.787282105115894
$ws2.Cells.Item(7,8).Value = 0.144097889752298
$ws2.Cells.Item(8,1).Value = 'AVGEXP1W_FOOD'
$ws2.Cells.Item(8,2).Value = 0.00227573677433688
$ws2.Cells.Item(8,3).Value = 0.0675031973532076
$ws2.Cells.Item(8,4).Value = 0.936090234045014
$ws2.Cells.Item(8,5).Value = -0.0116966668774625
$ws2.Cells.Item(8,6).Value = 0.0353391438933072
$ws2.Cells.Item(8,7).Value = 0.102811156286432
$ws2.Cells.Item(8,8).Value = 0.014725397531182
$ws2.Cells.Item(9,1).Value = 'AVGEXP1Y_TOTAL'
$ws2.Cells.Item(9,2).Value = 0.031985748372781
$ws2.Cells.Item(9,3).Value = 0.0370704985268565
$ws2.Cells.Item(9,4).Value = 0.912674793658593
$ws2.Cells.Item(9,5).Value = 0.1525899813309
$ws2.Cells.Item(9,6).Value = 0.0582794745938844
$ws2.Cells.Item(9,7).Value = 0.197210696644722
$ws2.Cells.Item(9,8).Value = -0.00494057001429787
$ws2.Cells.Item(10,1).Value = 'CHILDCARE_RATIO'
$ws2.Cells.Item(10,2).Value = 0.955148722884666
$ws2.Cells.Item(10,3).Value = -0.107951608566766
$ws2.Cells.Item(10,4).Value = -0.0513618660769178
$ws2.Cells.Item(10,5).Value = 0.0988477681041689
$ws2.Cells.Item(10,6).Value = -0.0629413776774466
$ws2.Cells.Item(10,7).Value = 0.0201977938821569
$ws2.Cells.Item(10,8).Value = -0.0128076411765843
$ws2.Cells.Item(11,1).Value = 'CHILDCORESD_RATIO'
$ws2.Cells.Item(11,2).Value = 0.724296648523608
$ws2.Cells.Item(11,3).Value = -0.223727428297577
$ws2.Cells.Item(11,4).Value = 0.27993235114076
$ws2.Cells.Item(11,5).Value = -0.10258985720099
$ws2.Cells.Item(11,6).Value = -0.21595193209955
$ws2.Cells.Item(11,7).Value = -0.117890758278472
$ws2.Cells.Item(11,8).Value = 0.0216074125642645
$ws2.Cells.Item(12,1).Value = 'CHILDLVNEAR_RATIO'
$ws2.Cells.Item(12,2).Value = 0.948915846579949
$ws2.Cells.Item(12,3).Value = -0.139747210908952
$ws2.Cells.Item(12,4).Value = -0.0570955495889937
$ws2.Cells.Item(12,5).Value = 0.084751059175772
$ws2.Cells.Item(12,6).Value = -0.07844367712414
$ws2.Cells.Item(12,7).Value = 0.00730817980701672
$ws2.Cells.Item(12,8).Value = -0.00537835625436425
$ws2.Cells.Item(13,1).Value = 'TRANSCHILD_RATIO'
$ws2.Cells.Item(13,2).Value = -0.20215214060006
$ws2.Cells.Item(13,3).Value = 0.639065983783579
$ws2.Cells.Item(13,4).Value = 0.0789068265441356
$ws2.Cells.Item(13,5).Value = -0.241869439669255
$ws2.Cells.Item(13,6).Value = 0.264256328763045
$ws2.Cells.Item(13,7).Value = 0.0540489520722349
$ws2.Cells.Item(13,8).Value = 0.331949611169065
$ws2.Cells.Item(14,1).Value = 'WORK_RATIO'
$ws2.Cells.Item(14,2).Value = 0.0218916442688938
$ws2.Cells.Item(14,3).Value = -0.164950990902395
$ws2.Cells.Item(14,4).Value = -0.109493700455018
$ws2.Cells.Item(14,5).Value = -0.594590694266201
$ws2.Cells.Item(14,6).Value = 0.303677130427794
$ws2.Cells.Item(14,7).Value = -0.407088393031316
$ws2.Cells.Item(14,8).Value = 0.391703351277459
$ws2.Cells.Item(15,1).Value = 'JOBSTATUS_NAGE_RATIO'
$ws2.Cells.Item(15,2).Value = 0.138387401196175
$ws2.Cells.Item(15,3).Value = -0.139499436729902
$ws2.Cells.Item(15,4).Value = 0.0787819660918628
$ws2.Cells.Item(15,5).Value = 0.763958979290977
$ws2.Cells.Item(15,6).Value = 0.181999842230487
$ws2.Cells.Item(15,7).Value = 0.0211806240751871
$ws2.Cells.Item(15,8).Value = -0.108427153227375
$ws2.Cells.Item(16,1).Value = 'JOBSTATUS_NAGS_RATIO'
$ws2.Cells.Item(16,2).Value = -0.0906185106024216
$ws2.Cells.Item(16,3).Value = -0.222854052885237
$ws2.Cells.Item(16,4).Value = 0.0254798753576386
$ws2.Cells.Item(16,5).Value = 0.604080761605345
$ws2.Cells.Item(16,6).Value = -0.0895121511274372
$ws2.Cells.Item(16,7).Value = -0.356329096946963
$ws2.Cells.Item(16,8).Value = 0.425951113393035
$ws2.Cells.Item(17,1).Value = 'JOBSTATUS_UNEM_RATIO'
$ws2.Cells.Item(17,2).Value = -0.037087153970814
$ws2.Cells.Item(17,3).Value = -0.21120690286187
$ws2.Cells.Item(17,4).Value = 0.00182614128114163
$ws2.Cells.Item(17,5).Value = 0.0694103995295483
$ws2.Cells.Item(17,6).Value = -0.0347393100437785
$ws2.Cells.Item(17,7).Value = -0.133907085996241
$ws2.Cells.Item(17,8).Value = -0.783306632942868
$ws2.Cells.Item(18,1).Value = 'JOBSTATUS_NEWK_RATIO'
$ws2.Cells.Item(18,2).Value = -0.00975696020297295
$ws2.Cells.Item(18,3).Value = -0.321287575883903
$ws2.Cells.Item(18,4).Value = -0.142109232792578
$ws2.Cells.Item(18,5).Value = 0.0629282736389248
$ws2.Cells.Item(18,6).Value = -0.698505076857583
$ws2.Cells.Item(18,7).Value = 0.00124768421458535
$ws2.Cells.Item(18,8).Value = -0.18444909446362
